$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Updated publish date
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into the new "Jurisdiction" / "United States of America" row.
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was the second (duplicate) "Contact" row - remove it entirely,
# shifting everything below up by one.
$ws1.Rows.Item(11).Delete()

# --- Elements sheet ---
$ws2 = $wb.Worksheets.Item("Elements")

# The root Extension row's Short/Definition columns now reflect this
# specific extension's Title/Description instead of the generic text.
$ws2.Range("K2").Value = "Effective Period"
$ws2.Range("L2").Value = "Effective period"
